$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the date string for C2 and C3 (11Jul2020 -> 12Jul2020)
$ws.Range("C2").Value = "12Jul2020"
$ws.Range("C3").Value = "12Jul2020"

# Update numeric values for columns K through BH, rows 2 and 3 (both rows share identical new values)
$ws.Range("K2").Value = 2
$ws.Range("K3").Value = 2
$ws.Range("L2").Value = 8
$ws.Range("L3").Value = 8
$ws.Range("M2").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O2").Value = 8.5906248092651367
$ws.Range("O3").Value = 8.5906248092651367
$ws.Range("P2").Value = 14.050758361816406
$ws.Range("P3").Value = 14.050758361816406
$ws.Range("Q2").Value = 12.506132125854492
$ws.Range("Q3").Value = 12.506132125854492
$ws.Range("R2").Value = 14.704978942871094
$ws.Range("R3").Value = 14.704978942871094
$ws.Range("S2").Value = 15.732240676879883
$ws.Range("S3").Value = 15.732240676879883
$ws.Range("T2").Value = 35.372631072998047
$ws.Range("T3").Value = 35.372631072998047
$ws.Range("U2").Value = 0
$ws.Range("U3").Value = 0
$ws.Range("V2").Value = 1
$ws.Range("V3").Value = 1
$ws.Range("W2").Value = 3.9138941764831543
$ws.Range("W3").Value = 3.9138941764831543
$ws.Range("X2").Value = 13.375571250915527
$ws.Range("X3").Value = 13.375571250915527
$ws.Range("Y2").Value = 44.65960693359375
$ws.Range("Y3").Value = 44.65960693359375
$ws.Range("Z2").Value = 7.6776456832885742
$ws.Range("Z3").Value = 7.6776456832885742
$ws.Range("AA2").Value = 38.660953521728516
$ws.Range("AA3").Value = 38.660953521728516
$ws.Range("AB2").Value = 30.264581680297852
$ws.Range("AB3").Value = 30.264581680297852
$ws.Range("AC2").Value = 1
$ws.Range("AC3").Value = 1
$ws.Range("AD2").Value = 1
$ws.Range("AD3").Value = 1
$ws.Range("AE2").Value = 19.735282897949219
$ws.Range("AE3").Value = 19.735282897949219
$ws.Range("AF2").Value = 39.1041259765625
$ws.Range("AF3").Value = 39.1041259765625
$ws.Range("AG2").Value = 3.3548846244812012
$ws.Range("AG3").Value = 3.3548846244812012
$ws.Range("AH2").Value = 4.9948983192443848
$ws.Range("AH3").Value = 4.9948983192443848
$ws.Range("AI2").Value = 31.555620193481445
$ws.Range("AI3").Value = 31.555620193481445
$ws.Range("AJ2").Value = 8.9748449325561523
$ws.Range("AJ3").Value = 8.9748449325561523
$ws.Range("AK2").Value = 1
$ws.Range("AK3").Value = 1
$ws.Range("AL2").Value = 1
$ws.Range("AL3").Value = 1
$ws.Range("AM2").Value = 3.3343899250030518
$ws.Range("AM3").Value = 3.3343899250030518
$ws.Range("AN2").Value = 12.788989067077637
$ws.Range("AN3").Value = 12.788989067077637
$ws.Range("AO2").Value = 32.278564453125
$ws.Range("AO3").Value = 32.278564453125
$ws.Range("AP2").Value = 12.842391967773438
$ws.Range("AP3").Value = 12.842391967773438
$ws.Range("AQ2").Value = 7.2027759552001953
$ws.Range("AQ3").Value = 7.2027759552001953
$ws.Range("AR2").Value = 37.126541137695313
$ws.Range("AR3").Value = 37.126541137695313
$ws.Range("AS2").Value = 1
$ws.Range("AS3").Value = 1
$ws.Range("AT2").Value = 0
$ws.Range("AT3").Value = 0
$ws.Range("AU2").Value = 14.665376663208008
$ws.Range("AU3").Value = 14.665376663208008
$ws.Range("AV2").Value = 4.4169931411743164
$ws.Range("AV3").Value = 4.4169931411743164
$ws.Range("AW2").Value = 3.0003960132598877
$ws.Range("AW3").Value = 3.0003960132598877
$ws.Range("AX2").Value = 8.2621965408325195
$ws.Range("AX3").Value = 8.2621965408325195
$ws.Range("AY2").Value = 49.364387512207031
$ws.Range("AY3").Value = 49.364387512207031
$ws.Range("AZ2").Value = 46.937423706054688
$ws.Range("AZ3").Value = 46.937423706054688
$ws.Range("BA2").Value = 0
$ws.Range("BA3").Value = 0
$ws.Range("BB2").Value = 1
$ws.Range("BB3").Value = 1
$ws.Range("BC2").Value = 1.1367254257202148
$ws.Range("BC3").Value = 1.1367254257202148
$ws.Range("BD2").Value = 20.284408569335938
$ws.Range("BD3").Value = 20.284408569335938
$ws.Range("BE2").Value = 30.025539398193359
$ws.Range("BE3").Value = 30.025539398193359
$ws.Range("BF2").Value = 5.3607869148254395
$ws.Range("BF3").Value = 5.3607869148254395
$ws.Range("BG2").Value = 34.062877655029297
$ws.Range("BG3").Value = 34.062877655029297
$ws.Range("BH2").Value = 37.349773406982422
$ws.Range("BH3").Value = 37.349773406982422
